# AutoCommit_22 ноября 2023 г. 8:58:11_SibNout2023
# Recreate the "grade reset" edit: homework-tracking columns (C:X) are
# renamed to a simple Дз1..Дз6 header and every student's marks are wiped,
# except two students who get a fresh 5/5/5/5/5/5 in the new Дз1..Дз6 slots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 3: collapse the 11 paired date columns (5сен..14ноя) into
#     6 single "Дз#" columns C..H; I..X keep their style but lose the text.
$ws.Range("C3:X3").ClearContents()
$ws.Range("C3").Value = "Дз1"
$ws.Range("D3").Value = "Дз2"
$ws.Range("E3").Value = "Дз3"
$ws.Range("F3").Value = "Дз4"
$ws.Range("G3").Value = "Дз5"
$ws.Range("H3").Value = "Дз6"

# --- Wipe all per-student marks (C:X), the row SUM formulas (Z) and the
#     rounded-grade column (AB) for every student row (5..34).
$ws.Range("C5:X34").ClearContents()
$ws.Range("Z5:Z34").ClearContents()
$ws.Range("AB5:AB34").ClearContents()

# --- Re-enter fresh marks for the two students who kept their homework
#     scores (row 22 "Салчук Виктория", row 33 "Шибаева Анна").
$ws.Range("C22:H22").Value = 5
$ws.Range("C33:H33").Value = 5

# --- Footer "weight" row (row 35) also cleared along with the data block.
$ws.Range("C35:X35").ClearContents()

# --- Restore the view: the frozen pane stays at C5/row4-col2, but the
#     last selected cell moves from AB29 to K16.
$ws.Range("K16").Select()
